# Refresh the crypto price (column D) and 1h volume change (column E)
# values, matching the automated "Updated cryptos list ... with GitHub
# Actions" data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.599.72"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "'2.606.52"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'592.64"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "'154.45"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.547"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "'2.605.77"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.125"
$ws.Range("E10").Value = "  +8.99%  "
$ws.Range("D11").Value = "'0.160"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "'5.22"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "'0.353"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("E14").Value = "  -3.75%  "
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "'3.077.23"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "'67.572.18"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "'2.601.33"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'11.14"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").Value = "'363.22"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").Value = "'7.63"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("E23").Value = "  -4.82%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'67.36"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "'9.73"
$ws.Range("E26").Value = "  -8.24%  "
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("D29").Value = "'574.36"
$ws.Range("E29").Value = "  -6.52%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "'1.42"
$ws.Range("E31").Value = "  -4.06%  "
$ws.Range("D32").Value = "'7.91"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'1.52"
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("D37").Value = "'4.93"
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("D38").Value = "'158.35"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("D39").Value = "'19.30"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").Value = "'0.369"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "'5.30"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "'2.54"
$ws.Range("E43").Value = "  -6.89%  "
$ws.Range("D44").Value = "'41.16"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'16.40"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "'155.21"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "'0.0₆0286"
$ws.Range("E48").Value = "  -6.94%  "
$ws.Range("D49").Value = "'3.74"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "'0.625"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").Value = "'20.71"
$ws.Range("E51").Value = "  -4.15%  "

# Clear the "quote prefix" text style Excel applies when a value is
# entered with a leading apostrophe, so cell formatting/style indices
# stay exactly as they were before the refresh.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
